$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 660171.0600000001
$ws.Range("I98").Value = 862377.9
$ws.Range("J98").Value = 2999
$ws.Range("K98").Value = 862377.9
$ws.Range("L98").Value = 2999
$ws.Range("M98").Value = -860879.9
$ws.Range("N98").Value = -5995

$ws.Range("H122").Value = 660171.0600000001
$ws.Range("I122").Value = 862377.9
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 2587133.7
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -2584683.7
$ws.Range("N122").Value = -13897

$ws.Range("H125").Value = 8009087
$ws.Range("I125").Value = 803.55554
$ws.Range("J125").Value = 22423996
$ws.Range("K125").Value = 7231.99986
$ws.Range("L125").Value = 201815964
$ws.Range("M125").Value = -4771.99986
$ws.Range("N125").Value = -201820884

$ws.Range("H127").Value = 1278.8889
$ws.Range("I127").Value = 976.6667
$ws.Range("J127").Value = 1339.3334
$ws.Range("K127").Value = 2930.0001
$ws.Range("L127").Value = 4018.0002
$ws.Range("M127").Value = 2029.9999
$ws.Range("N127").Value = -13938.0002

$ws.Range("H129").Value = 1014.7111
$ws.Range("J129").Value = 1031.0227
$ws.Range("L129").Value = 3093.0681
$ws.Range("N129").Value = -13093.0681

$ws.Range("H131").Value = 12289.277
$ws.Range("I131").Value = 3200.6365
$ws.Range("J131").Value = 26571.428
$ws.Range("K131").Value = 9601.9095
$ws.Range("L131").Value = 79714.284
$ws.Range("M131").Value = -4561.9095
$ws.Range("N131").Value = -89794.284

$ws.Range("H132").Value = 25416.28
$ws.Range("I132").Value = 25416.28
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 76248.84
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -73718.84
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 62502932
$ws.Range("I137").Value = 111113490
$ws.Range("J137").Value = 3643.2856
$ws.Range("K137").Value = 333340470
$ws.Range("L137").Value = 10929.8568
$ws.Range("M137").Value = -333337920
$ws.Range("N137").Value = -16029.8568

$ws.Range("H138").Value = 3511.697
$ws.Range("I138").Value = 1396.7646
$ws.Range("J138").Value = 3950.1584
$ws.Range("K138").Value = 4190.293799999999
$ws.Range("L138").Value = 11850.4752
$ws.Range("M138").Value = 949.7062000000005
$ws.Range("N138").Value = -22130.4752

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 58718.11
$ws.Range("I32").Value = 11208.333
$ws.Range("J32").Value = 225002.33
$ws.Range("K32").Value = 11208.333
$ws.Range("L32").Value = 225002.33
$ws.Range("M32").Value = -10921.333
$ws.Range("N32").Value = -225576.33

$ws.Range("H46").Value = 3068.4
$ws.Range("I46").Value = 2569
$ws.Range("J46").Value = 3401.3333
$ws.Range("K46").Value = 2569
$ws.Range("L46").Value = 3401.3333
$ws.Range("M46").Value = -2250
$ws.Range("N46").Value = -4039.3333

$ws.Range("H61").Value = 2900
$ws.Range("I61").Value = 2279.6155
$ws.Range("J61").Value = 4916.25
$ws.Range("K61").Value = 2279.6155
$ws.Range("L61").Value = 4916.25
$ws.Range("M61").Value = -2067.6155
$ws.Range("N61").Value = -5340.25

$ws.Range("H74").Value = 4689.027
$ws.Range("I74").Value = 955.2
$ws.Range("J74").Value = 12467.833
$ws.Range("K74").Value = 955.2
$ws.Range("L74").Value = 12467.833
$ws.Range("M74").Value = -81.20000000000005
$ws.Range("N74").Value = -14215.833

$ws.Range("H77").Value = 4689.027
$ws.Range("I77").Value = 955.2
$ws.Range("J77").Value = 12467.833
$ws.Range("K77").Value = 4776
$ws.Range("L77").Value = 62339.165
$ws.Range("M77").Value = -408
$ws.Range("N77").Value = -71075.16500000001

$ws.Range("H132").Value = 3007.4707
$ws.Range("I132").Value = 2854.5715
$ws.Range("J132").Value = 3721
$ws.Range("K132").Value = 8563.7145
$ws.Range("L132").Value = 11163
$ws.Range("M132").Value = -6033.7145
$ws.Range("N132").Value = -16223

$ws.Range("H136").Value = 2900
$ws.Range("I136").Value = 2279.6155
$ws.Range("J136").Value = 4916.25
$ws.Range("K136").Value = 6838.8465
$ws.Range("L136").Value = 14748.75
$ws.Range("M136").Value = -4288.8465
$ws.Range("N136").Value = -19848.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6485.0415
$ws.Range("I31").Value = 3174.875
$ws.Range("J31").Value = 8140.125
$ws.Range("K31").Value = 3174.875
$ws.Range("L31").Value = 8140.125
$ws.Range("M31").Value = -2879.875
$ws.Range("N31").Value = -8730.125

$ws.Range("H34").Value = 6485.0415
$ws.Range("I34").Value = 3174.875
$ws.Range("J34").Value = 8140.125
$ws.Range("K34").Value = 3174.875
$ws.Range("L34").Value = 8140.125
$ws.Range("M34").Value = -2972.875
$ws.Range("N34").Value = -8544.125

$ws.Range("H99").Value = 47621616
$ws.Range("I99").Value = 2993
$ws.Range("J99").Value = 333333340
$ws.Range("K99").Value = 2993
$ws.Range("L99").Value = 333333340
$ws.Range("M99").Value = -1495
$ws.Range("N99").Value = -333336336

$ws.Range("H122").Value = 1701.0526
$ws.Range("I122").Value = 1778.3846
$ws.Range("J122").Value = 1533.5
$ws.Range("K122").Value = 5335.1538
$ws.Range("L122").Value = 4600.5
$ws.Range("M122").Value = -2885.1538
$ws.Range("N122").Value = -9500.5

$ws.Range("H126").Value = 47621616
$ws.Range("I126").Value = 2993
$ws.Range("J126").Value = 333333340
$ws.Range("K126").Value = 8979
$ws.Range("L126").Value = 1000000020
$ws.Range("M126").Value = -6509
$ws.Range("N126").Value = -1000004960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1222090
$ws.Range("I122").Value = 724.75
$ws.Range("J122").Value = 1509470.1
$ws.Range("K122").Value = 6522.75
$ws.Range("L122").Value = 13585230.9
$ws.Range("M122").Value = -4072.75
$ws.Range("N122").Value = -13590130.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 10531.714
$ws.Range("I102").Value = 10012
$ws.Range("J102").Value = 10739.6
$ws.Range("K102").Value = 10012
$ws.Range("L102").Value = 10739.6
$ws.Range("M102").Value = -8390
$ws.Range("N102").Value = -13983.6

$ws.Range("H122").Value = 4871.5
$ws.Range("I122").Value = 4348
$ws.Range("K122").Value = 13044
$ws.Range("M122").Value = -10594

$ws.Range("H126").Value = 2857
$ws.Range("I126").Value = 2785.3333
$ws.Range("J126").Value = 2888.8518
$ws.Range("K126").Value = 8355.999899999999
$ws.Range("L126").Value = 8666.555399999999
$ws.Range("M126").Value = -5885.999899999999
$ws.Range("N126").Value = -13606.5554

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3246.3215
$ws.Range("I7").Value = 2599.7
$ws.Range("J7").Value = 3605.5557
$ws.Range("K7").Value = 2599.7
$ws.Range("L7").Value = 3605.5557
$ws.Range("M7").Value = -2487.7
$ws.Range("N7").Value = -3829.5557

$ws.Range("H40").Value = 2948.1428
$ws.Range("I40").Value = 1915.8572
$ws.Range("J40").Value = 3464.2856
$ws.Range("K40").Value = 1915.8572
$ws.Range("L40").Value = 3464.2856
$ws.Range("M40").Value = -1779.8572
$ws.Range("N40").Value = -3736.2856

$ws.Range("H68").Value = 2100
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 2233.3333
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 2233.3333
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -3731.3333

$ws.Range("H71").Value = 2100
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 2233.3333
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 11166.6665
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -18654.6665

$ws.Range("H122").Value = 3319.5
$ws.Range("I122").Value = 2558
$ws.Range("J122").Value = 3573.3333
$ws.Range("K122").Value = 7674
$ws.Range("L122").Value = 10719.9999
$ws.Range("M122").Value = -5224
$ws.Range("N122").Value = -15619.9999

$ws.Range("H126").Value = 3246.3215
$ws.Range("I126").Value = 2599.7
$ws.Range("J126").Value = 3605.5557
$ws.Range("K126").Value = 7799.099999999999
$ws.Range("L126").Value = 10816.6671
$ws.Range("M126").Value = -5329.099999999999
$ws.Range("N126").Value = -15756.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 584.04346
$ws.Range("J126").Value = 1449.5
$ws.Range("L126").Value = 4348.5
$ws.Range("N126").Value = -9288.5
